# Insert a new slide "Step 2: Wrangling the data [iv]" right after the
# existing "Step 2: Wrangling the data [iii]" slide (position 19) and
# before "Step 3: Splitting the data" (which was position 20 and shifts
# down to position 21).

$p = $ppt.ActivePresentation

# Layout 2 == "Title and Content" (same layout used by the neighbouring
# slides, e.g. the "Step 3: Splitting the data" slide).
$new = $p.Slides.Add(20, 2)

$new.Shapes.Item(1).TextFrame.TextRange.Text = "Step 2: Wrangling the data [iv]"

$body = $new.Shapes.Item(2).TextFrame.TextRange
$body.Text = "We attempted scaling the data and setting variance thresholds of 0.03 and 0.05.`rBut these only worsened our final predictions.`rSo our final models did not take it into consideration."
